# Update countries & provincias Spain
# - Reorders a few countries in the "Pais" list (Egipto/Irak, the
#   Namibia/Chipre/Uruguay/Georgia block, and Groenlandia/Islas Malvinas)
#   by refreshing the country name + stats that land on each existing row.
# - Refreshes a batch of per-country case counters.
# - Bumps the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Country name swaps (row keeps its position, label + stats move) ----
$ws.Range("A26").Value = "Irak"
$ws.Range("A27").Value = "Egipto"

$ws.Range("A143").Value = "Namibia"
$ws.Range("A144").Value = "Republica de Chipre"
$ws.Range("A145").Value = "Uruguay"
$ws.Range("A146").Value = "Georgia"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# ---- Updated per-country counters (Casos totales, Nuevos casos, Casos
#      activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----
$updates = @{
    "B4" = 3618739;  "C4" = 1912;  "D4" = 1646683; "E4" = 1831871; "G4" = 41;  "H4" = 140185;
    "B5" = 1972072;  "C5" = 1163;                  "E5" = 529729;  "G5" = 45;  "H5" = 75568;

    "D19" = 186400;  "E19" = 5729;

    "B26" = 86148;   "C26" = 2281; "D26" = 54316;  "E26" = 28310;  "G26" = 90; "H26" = 3522;
    "B27" = 84843;                 "D27" = 26135;  "E27" = 54641;              "H27" = 4067;

    "B37" = 57668;   "C37" = 791;  "D37" = 47545;  "E37" = 9721;   "G37" = 3;  "H37" = 402;

    "E67" = 6059;    "G67" = 2;    "H67" = 73;

    "B86" = 7681;    "C86" = 270;  "D86" = 3534;   "E86" = 3907;   "G86" = 5;  "H86" = 240;

    "B94" = 6089;    "C94" = 484;  "D94" = 2951;   "E94" = 3085;   "G94" = 10; "H94" = 53;

    "B100" = 4039;   "C100" = 86;  "D100" = 2729;  "E100" = 1190;

    "B134" = 1327;   "C134" = 8;   "D134" = 1093;  "E134" = 184;

    "B143" = 1032;   "C143" = 72;  "D143" = 31;    "E143" = 999;               "H143" = 2;
    "B144" = 1025;                 "D144" = 839;   "E144" = 167;               "H144" = 19;
    "B145" = 1009;   "C145" = 0;   "D145" = 909;   "E145" = 69;                "H145" = 31;
    "B146" = 1006;   "C146" = 2;   "D146" = 883;   "E146" = 108;               "H146" = 15;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Timestamp footer ----
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 14:32"
